$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet is protected; unprotect to allow edits, then restore protection after.
$ws.Unprotect("lido")

# Update the confidentiality / as-of-date note (shared string used by A7)
$ws.Range("A7").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-22 for illustrative purposes only and are subject to change."
# Undo the automatic row-height bump that results from the cell text change,
# restoring row 7 back to the (implicit) default height.
$ws.Rows(7).AutoFit()

# Update the allocation / change figures
$ws.Range("D2").Value = 0.8412644833684334
$ws.Range("E2").Value = 0.000915870731388102
$ws.Range("D3").Value = 0.1587355166315666
$ws.Range("E3").Value = -0.002958579881656709
$ws.Range("E4").Value = 0.0003008578116630023

# Restore sheet protection (same state as before the edits)
$ws.Protect("lido")
